$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 8, shifting the existing rows 8-11 down to 9-12.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new metier data
# (POL, BALT, OTB_DEF_>0_0_0, n_count=5, KG_sum=1300, EUR_sum=0)
$ws.Range("A8").Value = "POL"
$ws.Range("B8").Value = "BALT"
$ws.Range("C8").Value = "OTB_DEF_>0_0_0"
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 1300
$ws.Range("F8").Value = 0

# Update the values for the PTB_DEF_105-115_1_120 row, now shifted to row 10
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 600
